$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old column D ("Terms Typically Offered"),
# shifting it to column G.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill "NA" for the new columns across all data rows (2-27).
$ws.Range("D2:F27").Value = "NA"
